# Replace embedded line-breaks in "Packaging" (and one "Manufacturer") cells
# with single spaces, joining the wrapped text onto one line.
# This mirrors the author's fix "Aventis Pasteur" / "10 x 1 dose vial" / etc.
# duplicate-with-newline strings being collapsed into the existing
# single-line shared-string entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value  = "10 dose vial 10 x 1 dose syringe"
$ws.Range("D8").Value  = "10 x 1 dose vial"
$ws.Range("D14").Value = "5 x 1 dose vials 5 x 1 dose syringe"
$ws.Range("D15").Value = "1 dose vial 1 dose TIP-LOK 5 x 1 dose TIP-LOK"
$ws.Range("D16").Value = "Singe dose vial"
$ws.Range("D17").Value = "10 x 1dose vial 5 x 1 dose TIP-LOK"
$ws.Range("D25").Value = "10 x 1 dose vials 10 x 3 dose vial 5 x 1 dose syringe"
$ws.Range("D26").Value = "1 x 1 dose vial 5 x 1 dose Tiplok 25 x 1 dose Tiplok"
$ws.Range("H29").Value = "Aventis Pasteur"
$ws.Range("H30").Value = "Aventis Pasteur"
$ws.Range("D30").Value = "10 dose vials"
$ws.Range("D31").Value = "10 dose vial"
